# Gestion auto liste comp
# Insert a new column "Date confirm. Liste comp." (with its value-expression
# companion row) right before the existing "Lieu Préselection" column (W),
# shifting every column from W onward one position to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column at position 23 (W) - everything from W onward shifts right.
$ws.Columns.Item(23).Insert()

# Populate the header (row 1) and the templating-expression row (row 2)
# for the freshly inserted column.
$ws.Range("W1").Value = "Date confirm. Liste comp."
$ws.Range("W2").Value = '${form.datConfirmListCompFormStr}'

# Match the authored column width for the new column (~27 characters, best-fit).
$ws.Columns.Item(23).ColumnWidth = 26.166666666666668

# Reflect the author's final active selection on the frozen (bottom-left) pane.
$ws.Range("W3").Select()
